$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2; existing rows 2-11 shift down to 3-12.
$ws.Rows(2).Insert()

# The inserted row inherits formatting from the row above (header row).
# Reset it so it matches the plain style used by the other data rows.
$ws.Range("A2:T2").ClearFormats()

# Column D uses a date-number style (copy it from the row that was
# pushed down, which still carries the correct style index).
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new week's record (same market/product metadata as every
# other row, new variety/quality/volume/price data).
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44552
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 15500
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15750
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1050
$ws.Range("T2").Value = 15
